$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coerce the price/volume columns to Text before writing so Excel
# keeps values like "1.000" / "242.96" / "0.3166" as literal strings
# instead of re-parsing them as numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.823.02"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "1.893.37"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "0.7944"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").Value = "242.96"
$ws.Range("E6").Value = "  +0.56%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "0.3166"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "25.43"
$ws.Range("E9").Value = "  -4.10%  "
$ws.Range("D10").Value = "0.07052"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "0.08078"
$ws.Range("E11").Value = "  +0.59%  "
$ws.Range("D12").Value = "0.7682"
$ws.Range("E12").Value = "  +2.79%  "
$ws.Range("D13").Value = "1.909.62"
$ws.Range("E13").Value = "  +0.96%  "
$ws.Range("D14").Value = "5.351"
$ws.Range("E14").Value = "  +2.71%  "
$ws.Range("E15").Value = "  +0.30%  "
$ws.Range("D16").Value = "29.845.32"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("D17").Value = "6.014"
$ws.Range("E17").Value = "  +2.03%  "
$ws.Range("D18").Value = "13.87"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("D19").Value = "244.55"
$ws.Range("E19").Value = "  -0.15%  "
$ws.Range("D20").Value = "0.000007711"
$ws.Range("E20").Value = "  -0.63%  "
$ws.Range("D21").Value = "8.369"
$ws.Range("E21").Value = "  +20.88%  "
$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").Value = "2.145.71"
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("D25").Value = "0.1642"
$ws.Range("E25").Value = "  +5.08%  "
$ws.Range("D26").Value = "9.354"
$ws.Range("E26").Value = "  +1.73%  "
$ws.Range("D27").Value = "166.08"
$ws.Range("E27").Value = "  +0.06%  "
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "1.400"
$ws.Range("E30").Value = "  +2.50%  "
$ws.Range("D31").Value = "1.540"
$ws.Range("E31").Value = "  +1.30%  "
$ws.Range("E32").Value = "  +3.88%  "
$ws.Range("D33").Value = "0.05688"
$ws.Range("E33").Value = "  +0.82%  "
$ws.Range("D34").Value = "4.041"
$ws.Range("E34").Value = "  -0.80%  "
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").Value = "0.7382"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("D37").Value = "0.9985"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "2.628"
$ws.Range("E38").Value = "  -3.49%  "
$ws.Range("D39").Value = "0.01909"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "0.4409"
$ws.Range("E41").Value = "  -0.36%  "
$ws.Range("D42").Value = "72.48"
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").Value = "5.812"
$ws.Range("E43").Value = "  -2.39%  "
$ws.Range("D44").Value = "0.8412"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "1.034.33"
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("D47").Value = "103.11"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("D48").Value = "1.870"
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").Value = "9.970"
$ws.Range("E49").Value = "  +2.68%  "
$ws.Range("D50").Value = "7.424"
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("D51").Value = "2.041.97"
$ws.Range("E51").Value = "  +0.07%  "

# Restore the original (default) cell formatting/style so the saved
# cells look exactly like the rest of the untouched text cells.
$dataRange.NumberFormat = "General"
$dataRange.Style = "Normal"
